# feat: add 2022-Q3 data
# Insert a new "2022-Q3" sheet (a copy of the "2022-Q2" fund-holdings sheet,
# positioned right after "总计") with the new quarter's fund data, and
# update the "总计" summary sheet with a new leading row for 2022-Q3
# (pushing the existing quarterly summary rows down by one row).

$wb = $excel.ActiveWorkbook

# --- 1. Build the new "2022-Q3" worksheet -------------------------------
$zj  = $wb.Worksheets.Item(1)          # "总计" (summary) sheet, stays first
$src = $wb.Worksheets.Item("2022-Q2")  # template sheet to copy structure from

# Copy "2022-Q2" so it lands right after "总计"; this gives the new sheet
# the same headers/column layout/styles as the other quarterly fund sheets.
$src.Copy($null, $zj)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Overwrite the single data row with the 2022-Q3 fund holding.
# Fund code / numeric-looking figures are stored as text in the source
# data, so force text via a leading apostrophe to avoid numeric coercion,
# then strip the resulting quote-prefix style so the cell keeps plain
# "Normal" styling (matching every other text cell in this table).
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'009188"
$q3.Range("C2").Value = "鹏华股息精选混合"
$q3.Range("D2").Value = "'0.62"
$q3.Range("E2").Value = "'86.55"
$q3.Range("F2").Value = "'1.76"
$q3.Range("G2").Value = "'0.0109"
$q3.Range("H2").Value = 7
$q3.Range("B2:G2").Style = "Normal"

# --- 2. Update the "总计" summary sheet ----------------------------------
# Push the existing quarterly rows (B2:D7) down one row to B3:D8, keeping
# column A's running index (0,1,2,...) untouched since it is already
# sequential by row position.
$zj.Range("B2:D7").Copy()
$zj.Range("B3").PasteSpecial()

# New first data row: 2022-Q3 summary figures.
$zj.Range("B2").Value = "2022-Q3"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.01

# New last row (row 8) needs its index cell created; copy format+value from
# the row above then overwrite with the correct running index (6).
$zj.Range("A7").Copy($zj.Range("A8"))
$zj.Range("A8").Value = 6

# --- 3. Restore the originally-active tab --------------------------------
# Copying a sheet makes the new copy the active/selected tab; the workbook
# originally had its last sheet ("2021-Q1") selected, so put focus back.
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Select()
